# Marksheet update: correct marking scheme changed (per-correct-answer
# points raised from 3 to 5), which ripples into the "Total" row and the
# Corr/total marks summary shown in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row (row 11): points awarded per correct answer: 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row (row 12): total correct marks: 15 right * 5 pts = 75
$ws.Range("B12").Value = 75

# Corr/total marks summary: 75 correct out of new max of 140 (28 * 5)
$ws.Range("E12").Value = "75/140"
